$d = $word.ActiveDocument

$replacements = @(
    @("996×5=4980", "426×2=852"),
    @("202×3=606", "386×9=3474"),
    @("226×3=678", "309×4=1236"),
    @("453×4=1812", "615×7=4305"),
    @("948×4=3792", "525×5=2625"),
    @("135×7=945", "665×9=5985"),
    @("700×8=5600", "389×9=3501"),
    @("368×6=2208", "345×8=2760"),
    @("639×6=3834", "420×2=840"),
    @("997×7=6979", "171×7=1197"),
    @("215×5=1075", "716×7=5012"),
    @("875×7=6125", "474×2=948"),
    @("207×3=621", "343×6=2058"),
    @("428×8=3424", "614×2=1228"),
    @("524×7=3668", "824×4=3296"),
    @("578×3=1734", "150×2=300"),
    @("304×5=1520", "643×9=5787"),
    @("637×9=5733", "849×2=1698"),
    @("108×9=972", "161×2=322"),
    @("281×5=1405", "675×4=2700"),
    @("959×5=4795", "417×3=1251"),
    @("494×9=4446", "543×5=2715"),
    @("585×9=5265", "765×7=5355"),
    @("662×9=5958", "447×7=3129"),
    @("492×6=2952", "991×8=7928")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
